$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.260.60'
$ws.Range("E2").Value = '  -1.72%  '
$ws.Range("D3").Value = '1.583.64'
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("E4").Value = '  -0.29%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '209.77'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("E6").Value = '  -1.20%  '
$ws.Range("E8").Value = '  -1.15%  '
$ws.Range("E9").Value = '  -0.10%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '19.61'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -0.21%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0846'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("D12").Value = '1.806.07'
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("D13").Value = '1.580.80'
$ws.Range("E13").Value = '  -2.18%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '4.03'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("E15").Value = '  -1.04%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '64.68'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -0.73%  '
$ws.Range("D17").Value = '26.262.25'
$ws.Range("E17").Value = '  -1.62%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '7.23'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("E20").Value = '  -0.24%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '206.86'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -1.71%  '
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("E23").Value = '  -3.62%  '
$ws.Range("E24").Value = '  -1.21%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '144.55'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("E28").Value = '  -0.75%  '
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("E30").Value = '  -1.66%  '
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("E32").Value = '  -0.82%  '
$ws.Range("E33").Value = '  -0.75%  '
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.26'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +7.41%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '1.286.40'
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("E36").Value = '  -0.18%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.608'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("E39").Value = '  -1.51%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.817'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.61%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '5.51'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +1.81%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.770'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -1.28%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '62.29'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("D45").Value = '1.718.86'
$ws.Range("E45").Value = '  -1.20%  '
$ws.Range("E46").Value = '  -1.99%  '
$ws.Range("E47").Value = '  -0.39%  '
$ws.Range("E48").Value = '  +0.52%  '
$ws.Range("E49").Value = '  -1.60%  '
$ws.Range("E50").Value = '  -0.15%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '7.41'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.42%  '
